$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($addr in @("E34","F34")) {
    $r = $ws.Range($addr)
    $r.Borders.Item(8).LineStyle = -4142
    $r.Borders.Item(9).LineStyle = -4142
}
